$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.057915210723877
$ws.Range("B1").Value = 1.359564781188965
$ws.Range("C1").Value = 1.183177947998047
$ws.Range("D1").Value = 1.246571898460388
$ws.Range("E1").Value = 1.286367774009705
